# daily auto push: 2026-02-24 14:15 UTC
# A new observation row for 2026/02/24 (Tue, hour=20) needs to be inserted
# immediately above the existing 2026/12/29 block (row 866), shifting every
# subsequent row down by one and growing the sheet from A1:D907 to A1:D908.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 866; everything at/after 866 shifts down.
$ws.Rows(866).Insert()

# Populate the newly inserted row. The leading apostrophe on the date value
# forces it to stay a literal text string ("2026/02/24") instead of being
# auto-converted to an Excel date serial number, matching how every other
# date cell in column A is stored in this sheet.
$ws.Range("A866").Value = "'2026/02/24"
$ws.Range("B866").Value = "火"
$ws.Range("C866").Value = 20
$ws.Range("D866").Value = 201
